$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codeforiati:group-name" (col D) and "codeforiati:group-code" (col E)
# columns - including their header cells - need to swap places, for every
# used row on the sheet.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value = $eVal
    $eCell.Value = $dVal
}

# Row 91 (US-USAGOV) keeps its organisation name as "United States" in
# column B (it must not be affected by the group-name/group-code swap).
$ws.Cells.Item(91, 2).Value = "United States"
